$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace column A (filenames) with the location codes currently in column C
for ($r = 2; $r -le 68; $r++) {
    $loc = $ws.Cells.Item($r, 3).Text
    $ws.Cells.Item($r, 1).Value = $loc
}

# Update headers
$ws.Range("A1").Value = "Loc"
$ws.Range("B1").Value = "P_max"

# Remove the now-redundant Electrode Locations column (column C)
$ws.Range("C:C").Delete()
